# Feria Lagunitas de Puerto Montt - Ciboulette
# Insert one new daily-price record at row 284 (pushing the existing
# rows 284-340 down to 285-341), matching the weekly logic update.
#
# The new row reuses the (now shifted) original row's values, except for
# the date (column D) and the volume (column J), which are the only two
# cells that actually change for the newly-inserted record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 284:340 down to 285:341, opening up a blank row 284.
$ws.Rows.Item(284).Insert()

# Seed the new row 284 with a copy of what is now row 285 (the record
# that used to live at row 284 before the shift), so every column -
# including formatting such as the date style on column D - starts out
# identical to its neighbour.
$ws.Rows.Item(285).Copy()
$ws.Rows.Item(284).PasteSpecial()

# Now apply the two real edits for the new record.
$ws.Cells.Item(284, 4).Value = 45015   # D284: Fecha
$ws.Cells.Item(284, 10).Value = 120    # J284: Volumen
